$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" column values must stay text so exact formatting
# (leading/trailing zeros, etc.) is preserved.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.36"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.422"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05846"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.376"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.342"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8078"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9696"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1428"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07462"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03216"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03039"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.138"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09396"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001594"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04802"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005886"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006130"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004110"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009944"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.226"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3208"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1296"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03870"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1076"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002588"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003052"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006387"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005603"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1458"

# Other text field updates (coin names, links, labels).
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCXBestin24h"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
